$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Goal (per the OOXML diff):
#  1. The empty paragraph that immediately follows the paragraph ending in
#     "...The defendant having failed to appear the following is ordered:"
#     is removed, merging its (identical) paragraph properties away and
#     leaving a single paragraph.
#  2. A "_GoBack" bookmark (bookmarkStart/bookmarkEnd, id 0) is added at the
#     very end of that merged paragraph (right before its paragraph mark).
#  3. The "_GoBack" bookmark that used to sit by itself in the empty
#     paragraph right after "...blocking of motor vehicle registration or
#     transfer of registration.{% endif %}" is removed from there (it has
#     effectively moved to the location in step 2).
# ---------------------------------------------------------------------------

# --- Step 1: remove the pre-existing "_GoBack" bookmark -------------------
# (A document can only contain one bookmark with a given name, and Word
#  keeps a hidden "_GoBack" bookmark around marking the last edit location.
#  We delete it here and re-insert it at the correct spot below.)
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# --- Step 2: find the paragraph ending in "...is ordered:" ----------------
$findRange = $d.Content
$found = $findRange.Find.Execute("the following is ordered:", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$paraA = $findRange.Paragraphs(1)
$paraAEnd = $paraA.Range.End

# The very next paragraph is the empty paragraph that must be merged away:
# delete the paragraph mark that separates the two paragraphs.
$mergeRange = $d.Range($paraAEnd - 1, $paraAEnd)
$mergeRange.Delete()

# --- Step 3: insert the "_GoBack" bookmark at the end of the (now merged) --
# --- paragraph, right before its paragraph mark ----------------------------
# A collapsed Range exactly at a paragraph-end boundary anchors incorrectly
# in this runtime, so we use a tiny temporary placeholder character to build
# a non-collapsed Range, anchor the bookmark to it, then remove the
# placeholder again (the bookmark stays put).
$paraA = $findRange.Paragraphs(1)
$endPos = $paraA.Range.End
$placeholderPos = $d.Range($endPos - 1, $endPos - 1)
$placeholderPos.InsertAfter("X")

$paraA = $findRange.Paragraphs(1)
$endPos = $paraA.Range.End
$placeholderRange = $d.Range($endPos - 2, $endPos - 1)
$d.Bookmarks.Add("_GoBack", $placeholderRange)

$placeholderRange = $d.Range($endPos - 2, $endPos - 1)
$placeholderRange.Delete()
